$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("B3").Value = "L Totals"
$ws.Range("C3").Value = " L Minutes"

$ws.Range("B39").Value = "R Totals"
$ws.Range("C39").Value = "R Minutes"

# --- New column C values for L block (rows 4-38) ---
$ws.Range("C4").Value = 0.283
$ws.Range("C5").Value = 0.883
$ws.Range("C6").Value = 2.026
$ws.Range("C7").Value = 2.696
$ws.Range("C8").Value = 3.609
$ws.Range("C9").Value = 4.574
$ws.Range("C10").Value = 5.281
$ws.Range("C11").Value = 5.696
$ws.Range("C12").Value = 5.737
$ws.Range("C13").Value = 5.905
$ws.Range("C14").Value = 6.474
$ws.Range("C15").Value = 7.157
$ws.Range("C16").Value = 7.268
$ws.Range("C17").Value = 7.346
$ws.Range("C18").Value = 7.444
$ws.Range("C19").Value = 8.132999999999999
$ws.Range("C20").Value = 10.831
$ws.Range("C21").Value = 11.467
$ws.Range("C22").Value = 12.632
$ws.Range("C23").Value = 12.735
$ws.Range("C24").Value = 13.417
$ws.Range("C25").Value = 15.738
$ws.Range("C26").Value = 53.07
$ws.Range("C27").Value = 67.125
$ws.Range("C28").Value = 68.01600000000001
$ws.Range("C29").Value = 69.649
$ws.Range("C30").Value = 96.018
$ws.Range("C31").Value = 96.553
$ws.Range("C32").Value = 97.113
$ws.Range("C33").Value = 98.248
$ws.Range("C34").Value = 98.798
$ws.Range("C35").Value = 99.974
$ws.Range("C36").Value = 110.653
$ws.Range("C37").Value = 111.559
$ws.Range("C38").Value = 112.346

# --- New column C values for R block (rows 40-46) ---
$ws.Range("C40").Value = 0.283
$ws.Range("C41").Value = 0.883
$ws.Range("C42").Value = 2.026
$ws.Range("C43").Value = 2.696
$ws.Range("C44").Value = 3.609
$ws.Range("C45").Value = 4.574
$ws.Range("C46").Value = 5.281

# --- New summary rows 55-58 ---
$ws.Range("A55").Value = "L ~ Minutes"
$ws.Range("A56").Value = "0.283, 1.166, 3.192, 5.888, 9.497, 14.07, 19.351, 25.047, 30.784, 36.689, 43.163, 50.32, 57.588, 64.934, 72.378, 80.511, 91.342, 102.809, 115.441, 128.176, 141.593, 157.331, 210.401, 277.526, 345.542, 415.191, 511.209, 607.762, 704.875, 803.123, 901.921, 1001.895, 1112.548, 1224.107, 1336.452"
$ws.Range("A57").Value = "R ~ Minutes"
$ws.Range("A58").Value = "5.507, 13.581, 25.286, 38.755, 138.597, 238.456, 351.784"
